$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.70"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05630"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.470"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8054"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.043"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07314"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03202"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02934"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09258"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001677"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.201"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04718"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005826"
$ws.Range("E17").Value = "16OneONEWorstin24h"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006412"
$ws.Range("B19").Value = "UpBots"
$ws.Range("C19").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.007509"
$ws.Range("E19").Value = "18UpBotsUBXTBestin24h"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001054"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004107"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001504"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.977"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "GateToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.387"
$ws.Range("E24").Value = "23GateTokenGT"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.123"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3267"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("B27").Value = "ProBitToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1314"
$ws.Range("E27").Value = "26ProBitTokenPROB"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04152"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006870"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003510"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1036"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009036"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005648"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6818"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.01822"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002106"
